$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 for the new IPO entry (신한제13호스팩, 2024-04-22)
$ws.Rows(2).Insert()

# The inserted row inherits the header row's formatting; strip it back to the
# plain (unstyled) look used by the rest of the data rows.
$ws.Range("A2:Q2").ClearFormats()

# Force text storage so date-looking strings ("2024-04-22", etc.) are kept as
# plain text (matching the rest of the sheet) instead of being converted to
# Excel date serial numbers.
$ws.Range("A2:Q2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-04-22"
$ws.Range("B2").Value = "신한제13호스팩"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 60
$ws.Range("E2").Value = "신한"
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-04-11"
$ws.Range("P2").Value = "2024-04-15"
$ws.Range("Q2").Value = 2250000

# Remove the text-number-format override again so the cells end up with no
# explicit style, same as every other data row in the sheet.
$ws.Range("A2:Q2").ClearFormats()

# Drop the two trailing rows (비엔케이제2호스팩 / 유진스팩10호) that were
# removed from the dataset.
$ws.Rows("13:14").Delete()
